# Fruta / hortaliza, semanal
# Insert a new daily record at row 34 (shifts the existing rows 34-81 down
# to 35-82, growing the used range from A1:T81 to A1:T82) and populate the
# new row with the latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(34).Insert()

$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44413
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100102
$ws.Range("H34").Value = "Cítricos"
$ws.Range("I34").Value = 100102004
$ws.Range("J34").Value = "Mandarina"
$ws.Range("K34").Value = "Clementina"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 400
$ws.Range("N34").Value = 6500
$ws.Range("O34").Value = 6500
$ws.Range("P34").Value = 6500
$ws.Range("Q34").Value = "`$/bandeja 10 kilos"
$ws.Range("R34").Value = "Provincia de Limarí"
$ws.Range("S34").Value = 650
$ws.Range("T34").Value = 10
